$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.895.79"
$ws.Range("E2").Value = "  -1.41%  "
$ws.Range("D3").Value = "2.444.40"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.67"
$ws.Range("E5").Value = "  +0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.52"
$ws.Range("E6").Value = "  -1.84%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.530"
$ws.Range("D9").Value = "2.438.26"
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.18"
$ws.Range("E12").Value = "  -0.52%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.339"
$ws.Range("E13").Value = "  -3.14%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.94"
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000172"
$ws.Range("E15").Value = "  -1.67%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "2.893.40"
$ws.Range("E16").Value = "  +0.76%  "
$ws.Range("D17").Value = "61.826.42"
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("D18").Value = "2.445.02"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.61"
$ws.Range("E19").Value = "  -4.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.20"
$ws.Range("E20").Value = "  +0.76%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "324.82"
$ws.Range("E21").Value = "  -1.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.05"
$ws.Range("E22").Value = "  -2.06%  "
$ws.Range("E23").Value = "  -0.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.90"
$ws.Range("E24").Value = "  -6.30%  "
$ws.Range("E25").Value = "  -1.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.17"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "584.61"
$ws.Range("E27").Value = "  -7.81%  "
$ws.Range("E28").Value = "  +0.31%  "
$ws.Range("D30").Value = "0.0₃0939"
$ws.Range("E30").Value = "  -3.00%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.87"
$ws.Range("E31").Value = "  -2.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.38"
$ws.Range("E32").Value = "  -4.52%  "
$ws.Range("E34").Value = "  -3.17%  "
$ws.Range("E35").Value = "  +0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.77"
$ws.Range("E36").Value = "  -4.67%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "152.85"
$ws.Range("E37").Value = "  +4.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.373"
$ws.Range("E38").Value = "  -0.82%  "
$ws.Range("E39").Value = "  -4.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.32"
$ws.Range("E40").Value = "  -1.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.16"
$ws.Range("E41").Value = "  -2.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "42.37"
$ws.Range("E42").Value = "  -0.35%  "
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("E44").Value = "  -4.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.39"
$ws.Range("E45").Value = "  -5.11%  "
$ws.Range("D46").Value = "0.0₆0275"
$ws.Range("E46").Value = "  +17.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "140.80"
$ws.Range("E47").Value = "  -2.96%  "
$ws.Range("E48").Value = "  -3.72%  "
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0512"
$ws.Range("E50").Value = "  -2.71%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.62"
$ws.Range("E51").Value = "  -0.82%  "
